# Updates the cryptos price table (Sheet1) with the latest scraped
# values: Price (column D) and Volume(1h) percentage (column E) for
# each coin row, plus the Chainlink / WrappedEther row swap (rows 13-14).
# NumberFormat is forced to "@" (Text) before each write so values that
# look numeric (prices, percentages) are stored as literal text, matching
# the source data (which uses inline strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.396.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.096.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5268"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4429"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.62"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09345"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.77"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.168.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.561"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.916"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.44"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001161"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06696"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.340"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.429.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.309"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.86"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.19"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.518"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.141"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.655"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.269"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.866"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06810"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7011"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.348"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2224"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6871"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.38"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.385"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +19.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.640"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.239"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.37%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.219"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.11%  "
